# Weekly update: insert two new price records (Ajo / Rosado / "guarda")
# for date 2022-02-25 (serial 44617) right before the existing row 168,
# shifting all subsequent rows down by two (168-196 -> 170-198).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 168 (pushes old 168.. down to 170..)
$ws.Rows.Item(168).EntireRow.Insert()
$ws.Rows.Item(168).EntireRow.Insert()

# New row 168: Ajo, Rosado, 1a (guarda)
$ws.Range("A168").Value = 9
$ws.Range("B168").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C168").Value = "Metropolitana"
$ws.Range("D168").Value = 44617
$ws.Range("E168").Value = 13
$ws.Range("F168").Value = 100112003
$ws.Range("G168").Value = "Ajo"
$ws.Range("H168").Value = "Rosado"
$ws.Range("I168").Value = "1a (guarda)"
$ws.Range("J168").Value = 61
$ws.Range("K168").Value = 8500
$ws.Range("L168").Value = 9000
$ws.Range("M168").Value = 8754
$ws.Range("N168").Value = "$/trenza 50 unidades"
$ws.Range("O168").Value = "Provincia de Talagante"
$ws.Range("P168").Value = 1751
$ws.Range("Q168").Value = 5
$ws.Range("R168").Value = "Hortaliza"

# New row 169: Ajo, Rosado, 2a (guarda)
$ws.Range("A169").Value = 9
$ws.Range("B169").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44617
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112003
$ws.Range("G169").Value = "Ajo"
$ws.Range("H169").Value = "Rosado"
$ws.Range("I169").Value = "2a (guarda)"
$ws.Range("J169").Value = 25
$ws.Range("K169").Value = 7000
$ws.Range("L169").Value = 7500
$ws.Range("M169").Value = 7260
$ws.Range("N169").Value = "$/trenza 50 unidades"
$ws.Range("O169").Value = "Provincia de Talagante"
$ws.Range("P169").Value = 1452
$ws.Range("Q169").Value = 5
$ws.Range("R169").Value = "Hortaliza"
